{"js": "// Update the worksheet date and every two-digit multiplication problem.\n// Each \"find\" string is unique in the document (verified against the\n// source OOXML), so a direct search + replace per pair is unambiguous.\nconst replacements = [\n  [\"2025-10-20 Monday\", \"2025-10-21 Tuesday\"],\n  [\"87\u00d711=\", \"17\u00d772=\"],\n  [\"40\u00d727=\", \"80\u00d799=\"],\n  [\"25\u00d763=\", \"58\u00d798=\"],\n  [\"59\u00d768=\", \"91\u00d727=\"],\n  [\"95\u00d784=\", \"40\u00d745=\"],\n  [\"61\u00d714=\", \"90\u00d730=\"],\n  [\"30\u00d743=\", \"35\u00d773=\"],\n  [\"78\u00d761=\", \"37\u00d736=\"],\n  [\"48\u00d775=\", \"82\u00d796=\"],\n  [\"90\u00d756=\", \"97\u00d717=\"],\n  [\"47\u00d728=\", \"33\u00d736=\"],\n  [\"38\u00d714=\", \"42\u00d793=\"],\n  [\"56\u00d727=\", \"69\u00d752=\"],\n  [\"53\u00d792=\", \"96\u00d795=\"],\n  [\"71\u00d763=\", \"90\u00d732=\"],\n  [\"47\u00d786=\", \"74\u00d730=\"],\n  [\"47\u00d785=\", \"17\u00d714=\"],\n  [\"28\u00d747=\", \"27\u00d799=\"],\n  [\"31\u00d784=\", \"39\u00d792=\"],\n  [\"21\u00d727=\", \"31\u00d747=\"],\n  [\"76\u00d737=\", \"90\u00d713=\"],\n  [\"26\u00d735=\", \"39\u00d770=\"],\n  [\"17\u00d760=\", \"53\u00d728=\"],\n  [\"75\u00d749=\", \"97\u00d784=\"],\n  [\"26\u00d753=\", \"47\u00d745=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit multiplication problem.\n# Each \"find\" string is unique in the document, so a direct Find/Replace\n# per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-20 Monday\", \"2025-10-21 Tuesday\"),\n    @(\"87\u00d711=\", \"17\u00d772=\"),\n    @(\"40\u00d727=\", \"80\u00d799=\"),\n    @(\"25\u00d763=\", \"58\u00d798=\"),\n    @(\"59\u00d768=\", \"91\u00d727=\"),\n    @(\"95\u00d784=\", \"40\u00d745=\"),\n    @(\"61\u00d714=\", \"90\u00d730=\"),\n    @(\"30\u00d743=\", \"35\u00d773=\"),\n    @(\"78\u00d761=\", \"37\u00d736=\"),\n    @(\"48\u00d775=\", \"82\u00d796=\"),\n    @(\"90\u00d756=\", \"97\u00d717=\"),\n    @(\"47\u00d728=\", \"33\u00d736=\"),\n    @(\"38\u00d714=\", \"42\u00d793=\"),\n    @(\"56\u00d727=\", \"69\u00d752=\"),\n    @(\"53\u00d792=\", \"96\u00d795=\"),\n    @(\"71\u00d763=\", \"90\u00d732=\"),\n    @(\"47\u00d786=\", \"74\u00d730=\"),\n    @(\"47\u00d785=\", \"17\u00d714=\"),\n    @(\"28\u00d747=\", \"27\u00d799=\"),\n    @(\"31\u00d784=\", \"39\u00d792=\"),\n    @(\"21\u00d727=\", \"31\u00d747=\"),\n    @(\"76\u00d737=\", \"90\u00d713=\"),\n    @(\"26\u00d735=\", \"39\u00d770=\"),\n    @(\"17\u00d760=\", \"53\u00d728=\"),\n    @(\"75\u00d749=\", \"97\u00d784=\"),\n    @(\"26\u00d753=\", \"47\u00d745=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
